$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair rows reflecting the new ECs sending-cluster group (3 target
# clusters x 3 sending clusters), replacing the previous FAPs/sCs-only data.
$rowsData = @(
    @{Row=2; A="ECs"; B="Inhba"; C="Tgfbr3"; D="ECs"; E=2; F=0.6666666666666666; G=3.675031333333333; H=11.025094; I=0.2032371147293133; J=0.2032371147293133; K=3; L=1; M=14.85604233333333; N=44.568127; O=0.09286934904108346; P=0.09286934904108346; Q=54.59642106432644; R=491.367789578938; S=0.01887449854589932; T=0.01887449854589932},
    @{Row=3; A="ECs"; B="Inhba"; C="Tgfbr3"; D="FAPs"; E=2; F=0.6666666666666666; G=3.675031333333333; H=11.025094; I=0.2032371147293133; J=0.2032371147293133; K=3; L=1; M=114.2734143333333; N=342.820243; O=0.7143556381787382; P=0.7143556381787382; Q=419.9583782419824; R=3779.625404177842; S=0.145183578794064; T=0.145183578794064},
    @{Row=4; A="ECs"; B="Inhba"; C="Tgfbr3"; D="sCs"; E=2; F=0.6666666666666666; G=3.675031333333333; H=11.025094; I=0.2032371147293133; J=0.2032371147293133; K=3; L=1; M=30.83766366666667; N=92.512991; O=0.1927750127801784; P=0.1927750127801784; Q=113.3293802217949; R=1019.964421996154; S=0.03917903738934995; T=0.03917903738934995},
    @{Row=5; A="FAPs"; B="Inhba"; C="Tgfbr3"; D="ECs"; E=3; F=1; G=10.108494; H=30.325482; I=0.5590213983169419; J=0.5590213983169419; K=3; L=1; M=14.85604233333333; N=44.568127; O=0.09286934904108346; P=0.09286934904108346; Q=150.172214790246; R=1351.549933112214; S=0.05191595336173062; T=0.05191595336173062},
    @{Row=6; A="FAPs"; B="Inhba"; C="Tgfbr3"; D="FAPs"; E=3; F=1; G=10.108494; H=30.325482; I=0.5590213983169419; J=0.5590213983169419; K=3; L=1; M=114.2734143333333; N=342.820243; O=0.7143556381787382; P=0.7143556381787382; Q=1155.132123148014; R=10396.18910833213; S=0.3993400877502696; T=0.3993400877502696},
    @{Row=7; A="FAPs"; B="Inhba"; C="Tgfbr3"; D="sCs"; E=3; F=1; G=10.108494; H=30.325482; I=0.5590213983169419; J=0.5590213983169419; K=3; L=1; M=30.83766366666667; N=92.512991; O=0.1927750127801784; P=0.1927750127801784; Q=311.722338148518; R=2805.501043336662; S=0.1077653572049417; T=0.1077653572049417},
    @{Row=8; A="sCs"; B="Inhba"; C="Tgfbr3"; D="ECs"; E=3; F=1; G=4.298956; H=12.896868; I=0.2377414869537448; J=0.2377414869537448; K=3; L=1; M=14.85604233333333; N=44.568127; O=0.09286934904108346; P=0.09286934904108346; Q=63.86547232513734; R=574.789250926236; S=0.02207889713345351; T=0.02207889713345351},
    @{Row=9; A="sCs"; B="Inhba"; C="Tgfbr3"; D="FAPs"; E=3; F=1; G=4.298956; H=12.896868; I=0.2377414869537448; J=0.2377414869537448; K=3; L=1; M=114.2734143333333; N=342.820243; O=0.7143556381787382; P=0.7143556381787382; Q=491.2563801887694; R=4421.307421698924; S=0.1698319716344045; T=0.1698319716344045},
    @{Row=10; A="sCs"; B="Inhba"; C="Tgfbr3"; D="sCs"; E=3; F=1; G=4.298956; H=12.896868; I=0.2377414869537448; J=0.2377414869537448; K=3; L=1; M=30.83766366666667; N=92.512991; O=0.1927750127801784; P=0.1927750127801784; Q=132.5697592457987; R=1193.127833212188; S=0.04583061818588676; T=0.04583061818588676}
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Range("A$r").Value = $rd.A
    $ws.Range("B$r").Value = $rd.B
    $ws.Range("C$r").Value = $rd.C
    $ws.Range("D$r").Value = $rd.D
    $ws.Range("E$r").Value = $rd.E
    $ws.Range("F$r").Value = $rd.F
    $ws.Range("G$r").Value = $rd.G
    $ws.Range("H$r").Value = $rd.H
    $ws.Range("I$r").Value = $rd.I
    $ws.Range("J$r").Value = $rd.J
    $ws.Range("K$r").Value = $rd.K
    $ws.Range("L$r").Value = $rd.L
    $ws.Range("M$r").Value = $rd.M
    $ws.Range("N$r").Value = $rd.N
    $ws.Range("O$r").Value = $rd.O
    $ws.Range("P$r").Value = $rd.P
    $ws.Range("Q$r").Value = $rd.Q
    $ws.Range("R$r").Value = $rd.R
    $ws.Range("S$r").Value = $rd.S
    $ws.Range("T$r").Value = $rd.T
}
